$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The target cell (B11) currently shows the shared-string "R40".
# It needs to become the literal text "1" (not a number) while keeping
# its existing style/format (General, border, etc.) untouched.
#
# Assigning a plain numeric-looking string via .Value/.Value2 makes Excel
# auto-convert it to a real number, and forcing NumberFormat="@" (or an
# apostrophe-quoted entry) creates a brand new style record. Neither of
# those matches the intended edit, so instead we write a text formula
# that evaluates to the string "1", then convert that formula result to
# a static value in place (Copy + PasteSpecial values-only). That keeps
# the original cell style and yields a genuine text cell.
$target = $ws.Range("B11")
$target.Formula = "=""1"""
$target.Copy()
$target.PasteSpecial(-4163)
$excel.CutCopyMode = $false
